$wb = $excel.ActiveWorkbook

# --- compare_models: update model names (A,B columns) per row due to reordering, and metrics (C-I) ---
$ws = $wb.Worksheets.Item("compare_models")
$ws.Range("A2").Value = 'et'
$ws.Range("B2").Value = 'Extra Trees Regressor'
$ws.Range("C2").Value = 5.3473
$ws.Range("D2").Value = 59.9795
$ws.Range("E2").Value = 7.672
$ws.Range("F2").Value = 0.9191
$ws.Range("G2").Value = 0.1365
$ws.Range("H2").Value = 0.1026
$ws.Range("I2").Value = 0.076

$ws.Range("A3").Value = 'gbr'
$ws.Range("B3").Value = 'Gradient Boosting Regressor'
$ws.Range("C3").Value = 6.3353
$ws.Range("D3").Value = 81.63509999999999
$ws.Range("E3").Value = 8.9533
$ws.Range("F3").Value = 0.8908
$ws.Range("G3").Value = 0.1476
$ws.Range("H3").Value = 0.1134
$ws.Range("I3").Value = 0.05

$ws.Range("A4").Value = 'rf'
$ws.Range("B4").Value = 'Random Forest Regressor'
$ws.Range("C4").Value = 6.3026
$ws.Range("D4").Value = 82.40779999999999
$ws.Range("E4").Value = 9.026999999999999
$ws.Range("F4").Value = 0.8888
$ws.Range("G4").Value = 0.1485
$ws.Range("H4").Value = 0.1139
$ws.Range("I4").Value = 0.08599999999999999

$ws.Range("A5").Value = 'ridge'
$ws.Range("B5").Value = 'Ridge Regression'
$ws.Range("C5").Value = 7.2375
$ws.Range("D5").Value = 86.3532
$ws.Range("E5").Value = 9.2369
$ws.Range("F5").Value = 0.8827
$ws.Range("G5").Value = 0.1712
$ws.Range("H5").Value = 0.1381
$ws.Range("I5").Value = 0.018

$ws.Range("A6").Value = 'lightgbm'
$ws.Range("B6").Value = 'Light Gradient Boosting Machine'
$ws.Range("C6").Value = 6.6708
$ws.Range("D6").Value = 86.504
$ws.Range("E6").Value = 9.286799999999999
$ws.Range("F6").Value = 0.883
$ws.Range("G6").Value = 0.1578
$ws.Range("H6").Value = 0.123
$ws.Range("I6").Value = 0.034

$ws.Range("A7").Value = 'lr'
$ws.Range("B7").Value = 'Linear Regression'
$ws.Range("C7").Value = 7.4156
$ws.Range("D7").Value = 88.22190000000001
$ws.Range("E7").Value = 9.316800000000001
$ws.Range("F7").Value = 0.8799
$ws.Range("G7").Value = 0.1814
$ws.Range("H7").Value = 0.1448
$ws.Range("I7").Value = 1.256

$ws.Range("A8").Value = 'ada'
$ws.Range("B8").Value = 'AdaBoost Regressor'
$ws.Range("C8").Value = 7.0202
$ws.Range("D8").Value = 90.20359999999999
$ws.Range("E8").Value = 9.3918
$ws.Range("F8").Value = 0.8782
$ws.Range("G8").Value = 0.1664
$ws.Range("H8").Value = 0.1361
$ws.Range("I8").Value = 0.05

$ws.Range("A9").Value = 'br'
$ws.Range("B9").Value = 'Bayesian Ridge'
$ws.Range("C9").Value = 7.4111
$ws.Range("D9").Value = 90.4932
$ws.Range("E9").Value = 9.470499999999999
$ws.Range("F9").Value = 0.8774
$ws.Range("G9").Value = 0.1732
$ws.Range("H9").Value = 0.1403
$ws.Range("I9").Value = 0.018

$ws.Range("A10").Value = 'knn'
$ws.Range("B10").Value = 'K Neighbors Regressor'
$ws.Range("C10").Value = 6.5495
$ws.Range("D10").Value = 95.5809
$ws.Range("E10").Value = 9.6897
$ws.Range("F10").Value = 0.8718
$ws.Range("G10").Value = 0.1624
$ws.Range("H10").Value = 0.1224
$ws.Range("I10").Value = 0.022

$ws.Range("A11").Value = 'huber'
$ws.Range("B11").Value = 'Huber Regressor'
$ws.Range("C11").Value = 7.7872
$ws.Range("D11").Value = 97.6611
$ws.Range("E11").Value = 9.795
$ws.Range("F11").Value = 0.8667
$ws.Range("G11").Value = 0.193
$ws.Range("H11").Value = 0.1531
$ws.Range("I11").Value = 0.028

$ws.Range("A12").Value = 'lasso'
$ws.Range("B12").Value = 'Lasso Regression'
$ws.Range("C12").Value = 7.7385
$ws.Range("D12").Value = 107.1094
$ws.Range("E12").Value = 10.296
$ws.Range("F12").Value = 0.8554
$ws.Range("G12").Value = 0.1866
$ws.Range("H12").Value = 0.1471
$ws.Range("I12").Value = 0.59

$ws.Range("A13").Value = 'en'
$ws.Range("B13").Value = 'Elastic Net'
$ws.Range("C13").Value = 8.1068
$ws.Range("D13").Value = 115.9043
$ws.Range("E13").Value = 10.7155
$ws.Range("F13").Value = 0.8436
$ws.Range("G13").Value = 0.1872
$ws.Range("H13").Value = 0.1524
$ws.Range("I13").Value = 0.018

$ws.Range("A14").Value = 'omp'
$ws.Range("B14").Value = 'Orthogonal Matching Pursuit'
$ws.Range("C14").Value = 8.076599999999999
$ws.Range("D14").Value = 134.0358
$ws.Range("E14").Value = 11.4976
$ws.Range("F14").Value = 0.8179
$ws.Range("G14").Value = 0.1916
$ws.Range("H14").Value = 0.1508
$ws.Range("I14").Value = 0.02

$ws.Range("A15").Value = 'dt'
$ws.Range("B15").Value = 'Decision Tree Regressor'
$ws.Range("C15").Value = 8.599600000000001
$ws.Range("D15").Value = 161.026
$ws.Range("E15").Value = 12.5111
$ws.Range("F15").Value = 0.783
$ws.Range("G15").Value = 0.2069
$ws.Range("H15").Value = 0.1566
$ws.Range("I15").Value = 0.022

$ws.Range("A16").Value = 'par'
$ws.Range("B16").Value = 'Passive Aggressive Regressor'
$ws.Range("C16").Value = 10.2316
$ws.Range("D16").Value = 179.6481
$ws.Range("E16").Value = 13.1691
$ws.Range("F16").Value = 0.7618
$ws.Range("G16").Value = 0.2398
$ws.Range("H16").Value = 0.1939
$ws.Range("I16").Value = 0.018

$ws.Range("A17").Value = 'llar'
$ws.Range("B17").Value = 'Lasso Least Angle Regression'
$ws.Range("C17").Value = 13.9291
$ws.Range("D17").Value = 293.151
$ws.Range("E17").Value = 17.0799
$ws.Range("F17").Value = 0.6072
$ws.Range("G17").Value = 0.3047
$ws.Range("H17").Value = 0.2892
$ws.Range("I17").Value = 0.018

$ws.Range("A18").Value = 'lar'
$ws.Range("B18").Value = 'Least Angle Regression'
$ws.Range("C18").Value = 21.9477
$ws.Range("D18").Value = 1065.3555
$ws.Range("E18").Value = 27.3593
$ws.Range("F18").Value = -0.7437
$ws.Range("G18").Value = 0.4931
$ws.Range("H18").Value = 0.4504
$ws.Range("I18").Value = 0.03

$ws.Range("A19").Value = 'dummy'
$ws.Range("B19").Value = 'Dummy Regressor'
$ws.Range("C19").Value = 23.1892
$ws.Range("D19").Value = 750.8919
$ws.Range("E19").Value = 27.3639
$ws.Range("F19").Value = -0.0047
$ws.Range("G19").Value = 0.4806
$ws.Range("H19").Value = 0.4979
$ws.Range("I19").Value = 0.018

# --- tuned_1: update metrics ---
$ws = $wb.Worksheets.Item('tuned_1')
$ws.Range("B2").Value = 4.5824
$ws.Range("C2").Value = 53.5247
$ws.Range("D2").Value = 7.3161
$ws.Range("E2").Value = 0.9271
$ws.Range("F2").Value = 0.09959999999999999
$ws.Range("G2").Value = 0.0721
$ws.Range("B3").Value = 4.9402
$ws.Range("C3").Value = 63.9136
$ws.Range("D3").Value = 7.9946
$ws.Range("E3").Value = 0.9213
$ws.Range("F3").Value = 0.132
$ws.Range("G3").Value = 0.09089999999999999
$ws.Range("B4").Value = 5.5151
$ws.Range("C4").Value = 52.0426
$ws.Range("D4").Value = 7.2141
$ws.Range("E4").Value = 0.9167999999999999
$ws.Range("F4").Value = 0.117
$ws.Range("G4").Value = 0.09520000000000001
$ws.Range("B5").Value = 6.1661
$ws.Range("C5").Value = 67.7056
$ws.Range("D5").Value = 8.228300000000001
$ws.Range("E5").Value = 0.9204
$ws.Range("F5").Value = 0.14
$ws.Range("G5").Value = 0.1146
$ws.Range("B6").Value = 7.0187
$ws.Range("C6").Value = 99.0106
$ws.Range("D6").Value = 9.9504
$ws.Range("E6").Value = 0.8616
$ws.Range("F6").Value = 0.2276
$ws.Range("G6").Value = 0.1676
$ws.Range("B7").Value = 5.6445
$ws.Range("C7").Value = 67.2394
$ws.Range("D7").Value = 8.140700000000001
$ws.Range("E7").Value = 0.9094
$ws.Range("F7").Value = 0.1433
$ws.Range("G7").Value = 0.1081
$ws.Range("B8").Value = 0.8719
$ws.Range("C8").Value = 16.9689
$ws.Range("D8").Value = 0.9842
$ws.Range("E8").Value = 0.0241
$ws.Range("F8").Value = 0.0444
$ws.Range("G8").Value = 0.0327

# --- tuned_2: update metrics ---
$ws = $wb.Worksheets.Item('tuned_2')
$ws.Range("B2").Value = 5.1192
$ws.Range("C2").Value = 45.8964
$ws.Range("D2").Value = 6.7747
$ws.Range("E2").Value = 0.9375
$ws.Range("F2").Value = 0.0949
$ws.Range("G2").Value = 0.079
$ws.Range("B3").Value = 4.023
$ws.Range("C3").Value = 36.6777
$ws.Range("D3").Value = 6.0562
$ws.Range("E3").Value = 0.9548
$ws.Range("F3").Value = 0.0999
$ws.Range("G3").Value = 0.0718
$ws.Range("B4").Value = 5.5806
$ws.Range("C4").Value = 53.993
$ws.Range("D4").Value = 7.348
$ws.Range("E4").Value = 0.9136
$ws.Range("F4").Value = 0.1139
$ws.Range("G4").Value = 0.093
$ws.Range("B5").Value = 5.6551
$ws.Range("C5").Value = 60.3441
$ws.Range("D5").Value = 7.7681
$ws.Range("E5").Value = 0.929
$ws.Range("F5").Value = 0.1278
$ws.Range("G5").Value = 0.101
$ws.Range("B6").Value = 7.2627
$ws.Range("C6").Value = 95.93680000000001
$ws.Range("D6").Value = 9.794700000000001
$ws.Range("E6").Value = 0.8659
$ws.Range("F6").Value = 0.2087
$ws.Range("G6").Value = 0.1581
$ws.Range("B7").Value = 5.5281
$ws.Range("C7").Value = 58.5696
$ws.Range("D7").Value = 7.5484
$ws.Range("E7").Value = 0.9202
$ws.Range("F7").Value = 0.1291
$ws.Range("G7").Value = 0.1006
$ws.Range("B8").Value = 1.045
$ws.Range("C8").Value = 20.2991
$ws.Range("D8").Value = 1.2617
$ws.Range("E8").Value = 0.0302
$ws.Range("F8").Value = 0.0415
$ws.Range("G8").Value = 0.0305

# --- tuned_3: update metrics ---
$ws = $wb.Worksheets.Item('tuned_3')
$ws.Range("B2").Value = 4.6757
$ws.Range("C2").Value = 51.9334
$ws.Range("D2").Value = 7.2065
$ws.Range("E2").Value = 0.9293
$ws.Range("F2").Value = 0.0982
$ws.Range("G2").Value = 0.0721
$ws.Range("B3").Value = 5.2223
$ws.Range("C3").Value = 64.9123
$ws.Range("D3").Value = 8.056800000000001
$ws.Range("E3").Value = 0.92
$ws.Range("F3").Value = 0.1281
$ws.Range("G3").Value = 0.0922
$ws.Range("B4").Value = 5.6262
$ws.Range("C4").Value = 59.5939
$ws.Range("D4").Value = 7.7197
$ws.Range("E4").Value = 0.9046999999999999
$ws.Range("F4").Value = 0.1234
$ws.Range("G4").Value = 0.0964
$ws.Range("B5").Value = 6.3615
$ws.Range("C5").Value = 79.0622
$ws.Range("D5").Value = 8.8917
$ws.Range("E5").Value = 0.907
$ws.Range("F5").Value = 0.1461
$ws.Range("G5").Value = 0.1156
$ws.Range("B6").Value = 7.1906
$ws.Range("C6").Value = 100.9938
$ws.Range("D6").Value = 10.0496
$ws.Range("E6").Value = 0.8589
$ws.Range("F6").Value = 0.2233
$ws.Range("G6").Value = 0.1666
$ws.Range("B7").Value = 5.8153
$ws.Range("C7").Value = 71.2991
$ws.Range("D7").Value = 8.3849
$ws.Range("E7").Value = 0.904
$ws.Range("F7").Value = 0.1438
$ws.Range("G7").Value = 0.1086
$ws.Range("B8").Value = 0.8804
$ws.Range("C8").Value = 17.2909
$ws.Range("D8").Value = 0.9967
$ws.Range("E8").Value = 0.0243
$ws.Range("F8").Value = 0.0426
$ws.Range("G8").Value = 0.0321

# --- tuned_4: update metrics ---
$ws = $wb.Worksheets.Item('tuned_4')
$ws.Range("B2").Value = 7.1041
$ws.Range("C2").Value = 74.1725
$ws.Range("D2").Value = 8.612299999999999
$ws.Range("E2").Value = 0.899
$ws.Range("F2").Value = 0.147
$ws.Range("G2").Value = 0.1258
$ws.Range("B3").Value = 5.7651
$ws.Range("C3").Value = 61.5515
$ws.Range("D3").Value = 7.8455
$ws.Range("E3").Value = 0.9242
$ws.Range("F3").Value = 0.1673
$ws.Range("G3").Value = 0.1216
$ws.Range("B4").Value = 6.9176
$ws.Range("C4").Value = 79.0017
$ws.Range("D4").Value = 8.888299999999999
$ws.Range("E4").Value = 0.8736
$ws.Range("F4").Value = 0.1473
$ws.Range("G4").Value = 0.1195
$ws.Range("B5").Value = 7.3352
$ws.Range("C5").Value = 86.1867
$ws.Range("D5").Value = 9.2837
$ws.Range("E5").Value = 0.8986
$ws.Range("F5").Value = 0.1919
$ws.Range("G5").Value = 0.1463
$ws.Range("B6").Value = 8.8774
$ws.Range("C6").Value = 124.2084
$ws.Range("D6").Value = 11.1449
$ws.Range("E6").Value = 0.8264
$ws.Range("F6").Value = 0.2142
$ws.Range("G6").Value = 0.1796
$ws.Range("B7").Value = 7.1999
$ws.Range("C7").Value = 85.02419999999999
$ws.Range("D7").Value = 9.1549
$ws.Range("E7").Value = 0.8844
$ws.Range("F7").Value = 0.1735
$ws.Range("G7").Value = 0.1386
$ws.Range("B8").Value = 0.998
$ws.Range("C8").Value = 21.1743
$ws.Range("D8").Value = 1.1006
$ws.Range("E8").Value = 0.0331
$ws.Range("F8").Value = 0.0262
$ws.Range("G8").Value = 0.0226

# --- tuned_5: update metrics ---
$ws = $wb.Worksheets.Item('tuned_5')
$ws.Range("B2").Value = 5.8796
$ws.Range("C2").Value = 60.4954
$ws.Range("D2").Value = 7.7779
$ws.Range("E2").Value = 0.9176
$ws.Range("F2").Value = 0.12
$ws.Range("G2").Value = 0.096
$ws.Range("B3").Value = 6.0753
$ws.Range("C3").Value = 72.37009999999999
$ws.Range("D3").Value = 8.507099999999999
$ws.Range("E3").Value = 0.9108000000000001
$ws.Range("F3").Value = 0.1334
$ws.Range("G3").Value = 0.1018
$ws.Range("B4").Value = 5.5133
$ws.Range("C4").Value = 64.26439999999999
$ws.Range("D4").Value = 8.016500000000001
$ws.Range("E4").Value = 0.8972
$ws.Range("F4").Value = 0.1162
$ws.Range("G4").Value = 0.0896
$ws.Range("B5").Value = 5.8096
$ws.Range("C5").Value = 70.1314
$ws.Range("D5").Value = 8.3744
$ws.Range("E5").Value = 0.9175
$ws.Range("F5").Value = 0.1413
$ws.Range("G5").Value = 0.1059
$ws.Range("B6").Value = 7.0676
$ws.Range("C6").Value = 85.6944
$ws.Range("D6").Value = 9.257099999999999
$ws.Range("E6").Value = 0.8802
$ws.Range("F6").Value = 0.2068
$ws.Range("G6").Value = 0.1606
$ws.Range("B7").Value = 6.0691
$ws.Range("C7").Value = 70.5911
$ws.Range("D7").Value = 8.3866
$ws.Range("E7").Value = 0.9046999999999999
$ws.Range("F7").Value = 0.1435
$ws.Range("G7").Value = 0.1108
$ws.Range("B8").Value = 0.5309
$ws.Range("C8").Value = 8.642200000000001
$ws.Range("D8").Value = 0.506
$ws.Range("E8").Value = 0.0143
$ws.Range("F8").Value = 0.0329
$ws.Range("G8").Value = 0.0255

# --- blend_model: update metrics ---
$ws = $wb.Worksheets.Item('blend_model')
$ws.Range("B2").Value = 4.6959
$ws.Range("C2").Value = 43.396
$ws.Range("D2").Value = 6.5876
$ws.Range("E2").Value = 0.9409
$ws.Range("F2").Value = 0.09180000000000001
$ws.Range("G2").Value = 0.07340000000000001
$ws.Range("B3").Value = 4.7842
$ws.Range("C3").Value = 50.0878
$ws.Range("D3").Value = 7.0773
$ws.Range("E3").Value = 0.9383
$ws.Range("F3").Value = 0.119
$ws.Range("G3").Value = 0.08740000000000001
$ws.Range("B4").Value = 5.2297
$ws.Range("C4").Value = 53.0925
$ws.Range("D4").Value = 7.2865
$ws.Range("E4").Value = 0.9151
$ws.Range("F4").Value = 0.1093
$ws.Range("G4").Value = 0.0854
$ws.Range("B5").Value = 5.7064
$ws.Range("C5").Value = 62.2778
$ws.Range("D5").Value = 7.8916
$ws.Range("E5").Value = 0.9267
$ws.Range("F5").Value = 0.1316
$ws.Range("G5").Value = 0.103
$ws.Range("B6").Value = 7.1597
$ws.Range("C6").Value = 90.95959999999999
$ws.Range("D6").Value = 9.5373
$ws.Range("E6").Value = 0.8729
$ws.Range("F6").Value = 0.2102
$ws.Range("G6").Value = 0.1614
$ws.Range("B7").Value = 5.5152
$ws.Range("C7").Value = 59.9628
$ws.Range("D7").Value = 7.676
$ws.Range("E7").Value = 0.9188
$ws.Range("F7").Value = 0.1324
$ws.Range("G7").Value = 0.1021
$ws.Range("B8").Value = 0.8976
$ws.Range("C8").Value = 16.6454
$ws.Range("D8").Value = 1.0204
$ws.Range("E8").Value = 0.0247
$ws.Range("F8").Value = 0.041
$ws.Range("G8").Value = 0.0311

# --- stack_model: update metrics ---
$ws = $wb.Worksheets.Item('stack_model')
$ws.Range("B2").Value = 3.5128
$ws.Range("C2").Value = 24.9999
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0.966
$ws.Range("F2").Value = 0.0699
$ws.Range("G2").Value = 0.0553
$ws.Range("B3").Value = 4.5394
$ws.Range("C3").Value = 50.9388
$ws.Range("D3").Value = 7.1371
$ws.Range("E3").Value = 0.9372
$ws.Range("F3").Value = 0.117
$ws.Range("G3").Value = 0.0829
$ws.Range("B4").Value = 4.9059
$ws.Range("C4").Value = 41.9177
$ws.Range("D4").Value = 6.4744
$ws.Range("E4").Value = 0.9330000000000001
$ws.Range("F4").Value = 0.09959999999999999
$ws.Range("G4").Value = 0.07969999999999999
$ws.Range("B5").Value = 5.7632
$ws.Range("C5").Value = 63.5255
$ws.Range("D5").Value = 7.9703
$ws.Range("E5").Value = 0.9253
$ws.Range("F5").Value = 0.1278
$ws.Range("G5").Value = 0.104
$ws.Range("B6").Value = 6.4461
$ws.Range("C6").Value = 73.8369
$ws.Range("D6").Value = 8.5928
$ws.Range("E6").Value = 0.8968
$ws.Range("F6").Value = 0.189
$ws.Range("G6").Value = 0.1436
$ws.Range("B7").Value = 5.0335
$ws.Range("C7").Value = 51.0438
$ws.Range("D7").Value = 7.0349
$ws.Range("E7").Value = 0.9316
$ws.Range("F7").Value = 0.1206
$ws.Range("G7").Value = 0.0931
$ws.Range("B8").Value = 1.0101
$ws.Range("C8").Value = 16.9523
$ws.Range("D8").Value = 1.2464
$ws.Range("E8").Value = 0.0222
$ws.Range("F8").Value = 0.0394
$ws.Range("G8").Value = 0.0296

# --- pred_blend: update metrics (model stays Voting Regressor) ---
$ws = $wb.Worksheets.Item("pred_blend")
$ws.Range("C2").Value = 4.1655
$ws.Range("D2").Value = 32.5466
$ws.Range("E2").Value = 5.705
$ws.Range("F2").Value = 0.9347
$ws.Range("G2").Value = 0.09329999999999999
$ws.Range("H2").Value = 0.0696

# --- pred_stack: update metrics (model stays Stacking Regressor) ---
$ws = $wb.Worksheets.Item("pred_stack")
$ws.Range("C2").Value = 4.1441
$ws.Range("D2").Value = 31.1755
$ws.Range("E2").Value = 5.5835
$ws.Range("F2").Value = 0.9374
$ws.Range("G2").Value = 0.09030000000000001
$ws.Range("H2").Value = 0.06850000000000001

# --- pred_final: model changes from Stacking Regressor to Voting Regressor, update metrics ---
$ws = $wb.Worksheets.Item("pred_final")
$ws.Range("B2").Value = 'Voting Regressor'
$ws.Range("C2").Value = 1.4974
$ws.Range("D2").Value = 4.0136
$ws.Range("E2").Value = 2.0034
$ws.Range("F2").Value = 0.9941
$ws.Range("G2").Value = 0.0385
$ws.Range("H2").Value = 0.0273
